# Change the date placeholder text on every slide from "21/04/2020" to
# "28/05/2020". In the canonical OOXML the new date is split across five
# runs with alternating languages (en-CH / en-US / en-CH / en-US / en-CH),
# exactly as PowerPoint does when the field is retyped character-by-character
# with autocorrect toggling the run language:
#   "2" (en-CH) + "8" (en-US) + "/0" (en-CH) + "5" (en-US) + "/2020" (en-CH)
#
# TextRange.LanguageID always stamps the *first* run of the text frame in
# this host, so the runs are built back-to-front with InsertBefore: each
# newly inserted chunk becomes run #0 right when its language is set, which
# pushes the previously-tagged runs one slot to the right without touching
# their already-correct language.

$p = $ppt.ActivePresentation

function Set-SlideDate($slideIndex) {
    $slide = $p.Slides.Item($slideIndex)
    $shape = $slide.Shapes.Item("Date Placeholder 3")
    $tr = $shape.TextFrame.TextRange

    $tr.Text = "/2020"
    $tr.LanguageID = "en-CH"

    $r5 = $tr.InsertBefore("5")
    $r5.LanguageID = "en-US"

    $r0slash = $tr.InsertBefore("/0")
    $r0slash.LanguageID = "en-CH"

    $r8 = $tr.InsertBefore("8")
    $r8.LanguageID = "en-US"

    $r2 = $tr.InsertBefore("2")
    $r2.LanguageID = "en-CH"
}

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    Set-SlideDate $i
}
